# Update countries & provincias Spain
# Refresh COVID numbers for a set of countries on the "Pais" sheet and
# re-sort (by swapping the affected rows) so the "Casos totales" column
# (B) stays in descending order, matching the source data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Timestamp cell (A1) -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 16:57"

# --- Straight data refresh (country keeps its row / position) ----------
Set-Row 4   4501219 2876  2189592 2159217 0 90  152410   # Estados Unidos
Set-Row 6   1558447 26312 999318  524644  0 261 34485    # India
Set-Row 21  208339  388   192000  7132    0 0   9207     # Alemania
Set-Row 62  24892   372   14047   10287   0 7   558      # Serbia
Set-Row 63  23947   426   16785   6403    0 6   759      # Moldavia
Set-Row 92  7320    44    6103    1157    0 0   60       # Tayikistan
Set-Row 106 3738    29    1728    1907    0 0   103      # Malaui
Set-Row 114 3003    118   1591    1388    0 1   24       # Hong Kong
Set-Row 119 2588    33    2353    148     0 0   87       # Cuba
Set-Row 168 351     1     294     51      0 0   6        # Birmania
Set-Row 192 63      0     11      51      0 1   1        # Papua Nueva Guinea

# --- Re-sorted blocks ----------------------------------------------------
# Republica Dominicana overtakes Kuwait (rows 40/41 swap, RD gets new data,
# Kuwait keeps its previous totals but drops one place).
$ws.Cells.Item(40, 1).Value = "Republica Dominicana"
Set-Row 40 66182 1492 33947 31112 0 22 1123
$ws.Cells.Item(41, 1).Value = "Kuwait"
Set-Row 41 65903 754  56467 8992  0 2  444

# Azerbaiyan overtakes Japon (rows 58/59 swap).
$ws.Cells.Item(58, 1).Value = "Azerbaiyan"
Set-Row 58 31221 363 24495 6288 0 8 438
$ws.Cells.Item(59, 1).Value = "Japon"
Set-Row 59 30961 0   22811 7152 0 0 998

# Namibia jumps two places ahead of Guinea-Bisau and Ruanda (rows 128-130
# cascade down by one).
$ws.Cells.Item(128, 1).Value = "Namibia"
Set-Row 128 1986 69 104  1873 0 1 9
$ws.Cells.Item(129, 1).Value = "Guinea-Bisau"
Set-Row 129 1954 0  803  1125 0 0 26
$ws.Cells.Item(130, 1).Value = "Ruanda"
Set-Row 130 1926 0  1005 916  0 0 5
